$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 192
$ws.Range("F4").Value = 327
$ws.Range("F5").Value = 396
$ws.Range("F6").Value = 2175
$ws.Range("F7").Value = 1519
$ws.Range("F8").Value = 1333
$ws.Range("F9").Value = 3053
$ws.Range("F10").Value = 2037
$ws.Range("F11").Value = 1438
$ws.Range("F12").Value = 1766
$ws.Range("F13").Value = 76
$ws.Range("F14").Value = 692
$ws.Range("F15").Value = 619
$ws.Range("F17").Value = 2234
$ws.Range("F18").Value = 1255
$ws.Range("F19").Value = 227
$ws.Range("F20").Value = 2228
$ws.Range("F21").Value = 1856
$ws.Range("F22").Value = 673
$ws.Range("F23").Value = 5479
$ws.Range("F24").Value = 1041
$ws.Range("F25").Value = 95
$ws.Range("F26").Value = 87
$ws.Range("F27").Value = 1171
$ws.Range("F28").Value = 242
$ws.Range("F29").Value = 1058
$ws.Range("F30").Value = 542
$ws.Range("F31").Value = 115
$ws.Range("F32").Value = 260
$ws.Range("F33").Value = 1175
$ws.Range("F35").Value = 3502
$ws.Range("F36").Value = 602
$ws.Range("F37").Value = 1126
$ws.Range("F38").Value = 57
$ws.Range("F39").Value = 84
$ws.Range("F40").Value = 935
$ws.Range("F41").Value = 1229
$ws.Range("F43").Value = 26
$ws.Range("F44").Value = 804
$ws.Range("F45").Value = 1017
$ws.Range("F49").Value = 22

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 148659
$ws.Range("F12").Value = 413
$ws.Range("F16").Value = 324
$ws.Range("F18").Value = 121
$ws.Range("F19").Value = 101
$ws.Range("F21").Value = 100
$ws.Range("F26").Value = 534
$ws.Range("F27").Value = 180
$ws.Range("F28").Value = 326
$ws.Range("F31").Value = 53
$ws.Range("F32").Value = 53
$ws.Range("F36").Value = 8
$ws.Range("F37").Value = 126
$ws.Range("F40").Value = 18
$ws.Range("F43").Value = 136

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 3261
$ws.Range("F5").Value = 366
$ws.Range("F7").Value = 919
$ws.Range("F8").Value = 1404
$ws.Range("F9").Value = 747
$ws.Range("F10").Value = 348
$ws.Range("F11").Value = 2621
$ws.Range("F12").Value = 180
$ws.Range("F13").Value = 277
$ws.Range("F14").Value = 1002

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 366
$ws.Range("F3").Value = 919
$ws.Range("F4").Value = 747
$ws.Range("F5").Value = 192
$ws.Range("F6").Value = 327
$ws.Range("F7").Value = 348
$ws.Range("F8").Value = 2621
$ws.Range("F9").Value = 2175
$ws.Range("F10").Value = 1519
$ws.Range("F11").Value = 1333
$ws.Range("F12").Value = 3053
$ws.Range("F13").Value = 2037
$ws.Range("F14").Value = 1438
$ws.Range("F16").Value = 1766
$ws.Range("F17").Value = 692
$ws.Range("F18").Value = 619
$ws.Range("F19").Value = 324
$ws.Range("F21").Value = 2234
$ws.Range("F22").Value = 180
$ws.Range("F23").Value = 1255
$ws.Range("F24").Value = 227
$ws.Range("F25").Value = 2228
$ws.Range("F26").Value = 5479
$ws.Range("F27").Value = 277
$ws.Range("F28").Value = 1041
$ws.Range("F29").Value = 95
$ws.Range("F30").Value = 1002
$ws.Range("F31").Value = 1171
$ws.Range("F32").Value = 326
$ws.Range("F33").Value = 53
$ws.Range("F34").Value = 1058
$ws.Range("F35").Value = 542
$ws.Range("F36").Value = 115
$ws.Range("F37").Value = 260
$ws.Range("F39").Value = 3502
$ws.Range("F40").Value = 602
$ws.Range("F42").Value = 1126
$ws.Range("F43").Value = 84
$ws.Range("F44").Value = 935
$ws.Range("F45").Value = 1229
$ws.Range("F46").Value = 26
$ws.Range("F47").Value = 804
$ws.Range("F48").Value = 1017
$ws.Range("F49").Value = 136
$ws.Range("F50").Value = 136
